$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells that are being updated to remain text,
# so numeric-looking strings like "1.000" are not coerced into numbers.
$ws.Range("D2:D33").NumberFormat = "@"
$ws.Range("D35:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.244.00"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.857.77"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "330.85"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4539"
$ws.Range("E7").Value = "  -3.53%  "
$ws.Range("D8").Value = "0.3916"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "47.70"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "0.07806"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "0.9781"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "21.39"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "1.844.89"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "5.797"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "6.966"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "87.64"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").Value = "0.06529"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "0.00001016"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "17.00"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("D21").Value = "1.016"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").Value = "28.190.21"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "5.273"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").Value = "10.59"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "2.255"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").Value = "2.064.55"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "156.12"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("D28").Value = "19.12"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").Value = "2.034"
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").Value = "5.256"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("D31").Value = "116.26"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").Value = "0.09252"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").Value = "0.9335"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "1.377"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "5.180"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "0.06021"
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").Value = "0.02183"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").Value = "8.146"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").Value = "1.167"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "0.5636"
$ws.Range("E42").Value = "  -5.02%  "
$ws.Range("D43").Value = "9.968"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("D44").Value = "0.1785"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("D45").Value = "1.248"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "2.307"
$ws.Range("E46").Value = "  +25.82%  "
$ws.Range("D47").Value = "0.07202"
$ws.Range("E47").Value = "  +4.74%  "
$ws.Range("D48").Value = "0.5358"
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").Value = "11.69"
$ws.Range("E49").Value = "  -3.71%  "
$ws.Range("D50").Value = "1.866"
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("D51").Value = "109.43"
$ws.Range("E51").Value = "  -1.87%  "
